$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct cell updates reproducing the diff.
# For column D (Price) values that look numeric, force text type first
# via NumberFormat "@" then restore the default "Normal" style so the
# cell keeps its original (unstyled) appearance but the stored value stays a string.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.460.99"
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.590.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.60%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "654.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.95%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.48"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.407"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.72%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.10%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.586.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.47%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.24%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.201"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.97%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.262.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.55%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.327.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.95%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000256"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.92%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.583.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.23%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.90%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.88%  "

# Row 22
$ws.Range("E22").Value = "  +7.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.487"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.34%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "511.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.58%  "

# Row 25
$ws.Range("E25").Value = "  +7.03%  "

# Row 26
$ws.Range("E26").Value = "  +2.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.50%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.781.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.39%  "

# Row 30
$ws.Range("E30").Value = "  +18.54%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.40%  "

# Row 32
$ws.Range("E32").Value = "  -0.23%  "

# Row 33
$ws.Range("E33").Value = "  +4.18%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.97%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.176"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.31%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.04%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.559"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.35%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.25%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "572.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.33%  "

# Row 41
$ws.Range("E41").Value = "  -0.01%  "

# Row 42
$ws.Range("E42").Value = "  +1.80%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.928"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.11%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.05%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.24%  "

# Row 46
$ws.Range("E46").Value = "  +0.49%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0419"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.10%  "

# Row 48
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.81%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +31.71%  "

# Row 50
$ws.Range("E50").Value = "  +1.29%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.72%  "
